# FunctionMap.xlsx polishing pass: remove the flushing "UBIDs helper" row from the
# Detailed View pipeline table, and correct the ShpInf_01 flag description on the
# Shape Flags sheet (it now documents multipolygon handling instead of the old
# "multiple lines" wording).

$wb = $excel.ActiveWorkbook

# --- Detailed View: drop the helper.CREDA_Project("UBIDs", infile) row -------
$detailed = $wb.Worksheets.Item("Detailed View")
$detailed.Rows.Item(18).Select()
$detailed.Rows.Item(18).Delete()

# --- Shape Flags: update the ShpInf_01 explanation/notes ----------------------
$shapeFlags = $wb.Worksheets.Item("Shape Flags")
$shapeFlags.Cells.Item(3, 2).Value = "Multipolygon record"
$shapeFlags.Cells.Item(3, 3).Value = "Set whenever multiple ShapeIDZs are connected by a single ShapeID via multipolgon input"
$shapeFlags.Range("B4").Select()

# Shape Flags ends up as the active/selected sheet after the edits.
$shapeFlags.Activate()
